$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells that hold numeric-looking text retain their
# original plain-text representation (e.g. trailing zeros, scientific-looking
# decimals) instead of being auto-converted to numbers by Excel.

$ws.Range("D2").Value = "69.765.50"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "3.821.84"
$ws.Range("E3").Value = "  +1.92%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.41"
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.58"
$ws.Range("E6").Value = "  -1.53%  "

$ws.Range("D7").Value = "3.817.33"
$ws.Range("E7").Value = "  +1.82%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.45"
$ws.Range("E11").Value = "  +2.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.486"
$ws.Range("E12").Value = "  -0.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.10"
$ws.Range("E13").Value = "  -2.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000256"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").Value = "4.428.25"
$ws.Range("E15").Value = "  +1.33%  "

$ws.Range("D16").Value = "3.797.08"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").Value = "69.754.28"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.59"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("E19").Value = "  -3.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.67"
$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "511.38"
$ws.Range("E21").Value = "  +1.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.79"
$ws.Range("E22").Value = "  +2.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.739"
$ws.Range("E23").Value = "  +1.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.49"
$ws.Range("E24").Value = "  -0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.55"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000144"
$ws.Range("E26").Value = "  +4.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.88"
$ws.Range("E27").Value = "  -2.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.64"
$ws.Range("E28").Value = "  -4.32%  "

$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.53"
$ws.Range("E30").Value = "  +0.73%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.00"
$ws.Range("E31").Value = "  +3.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("E32").Value = "  +2.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.66"
$ws.Range("E33").Value = "  +1.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("E34").Value = "  -0.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.16"
$ws.Range("E37").Value = "  -0.67%  "

$ws.Range("E38").Value = "  +7.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "485.40"
$ws.Range("E39").Value = "  +13.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.341"
$ws.Range("E40").Value = "  +1.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.08"
$ws.Range("E41").Value = "  -1.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("E42").Value = "  +4.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.85"
$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.38"
$ws.Range("E44").Value = "  -2.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.61"
$ws.Range("E45").Value = "  -1.42%  "

$ws.Range("D46").Value = "2.953.34"
$ws.Range("E46").Value = "  -1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0364"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.56"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.00"
$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.48"
$ws.Range("E51").Value = "  -0.77%  "
